$d = $word.ActiveDocument

# --- Step 1: merge the first paragraph's runs and drop the proofErr spellcheck markers ---
# A find/replace across the whole paragraph (which spans the proofErr-wrapped run) collapses
# the three runs + proofErr tags into one clean run.
$d.Content.Find.Execute('Max Kelly, Nithin Perumal, Justin Hohl, Sam Weiskettal, Tim Smith', $true, $false, $false, $false, $false, $true, 1, $false, 'Max Kelly, Nithin Perumal, Justin Hohl, Sam Weiskettal, Tim Smith', 2) | Out-Null

# --- Step 2: rewrite the GUI-status paragraph (6) with the extended narrative, run-for-run ---
$guiPara = $d.Paragraphs.Item(6).Range
$guiPara.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t xml:space="preserve">The current state of the project is going a bit slower than planned.  Our GUI has been designed </w:t></w:r><w:r><w:t xml:space="preserve">and </w:t></w:r><w:r><w:t>operated by Tim Smith</w:t></w:r><w:r><w:t>.  He has been working on making the GUI user friendly and clickable so that you know which nodes you are on and characteristics about each node such as its x and y coordinates on the map.  The GUI also has several buttons in case you wish to select a particular node and </w:t></w:r><w:r><w:t xml:space="preserve">calculate a planned trip.  While the buttons themselves are currently inoperable, this is because the methods to perform their operations are largely unwritten at this current time. </w:t></w:r><w:r><w:t xml:space="preserve"> The nodes have an arraylist of neighbors that takes in the nodes </w:t></w:r><w:r><w:t>as they are created.  When they are clicked on, their color will change so that the user knows which nodes are selected and will be used to calculate the trip.  We have come up with a way to calculate the shortest distances with A*, but it has not been typed up yet in code.  Nithin has been working on text documents that are specifically written so that the program will load them and read them so that it will use what is written in the text documents to help program the nodes that will then be written onto the GUI.</w:t></w:r></w:p>')

# --- Step 3: append the new closing paragraph (tab + two sentences) ---
$guiPara = $d.Paragraphs.Item(6).Range
$guiPara.Collapse(0)
$guiPara.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(7).Range
$newPara.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:tab/><w:t>We realize that we are behind most of the other groups.  This project has been a bit of a struggle for us due to most of the group’s limited availability</w:t></w:r><w:r><w:t>.  We hope to get all of the requirements out of the way this upcoming week and begin work on some of our more interesting and fun features.  Currently, we still need to get our routes fully calculated and then finish our desired time calculations when traveling distance.  We also still need to implement the exact graph for our full trip planner.  Most of the routes that we have planned to implement will require a queue and a heap.  This will all be accomplished in the upcoming week</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>')

# --- Step 4: move the _GoBack bookmark to the very end of the (now last) paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$newPara = $d.Paragraphs.Item(7).Range
$endPos = $newPara.End - 1
$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output $d.Content.Text
